$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (vehicle_model): switch from select_by_text/"A4" to select_by_index/"3"
$ws.Range("D4").Value = "select_by_index"
$ws.Range("E4").Value = "3"

# Row 7 (vehicle_confirm): fix element value from "fake_vehicle_confirm" to "vehicle_confirm"
$ws.Range("C7").Value = "vehicle_confirm"

# Update the active selection to match the resaved file (cursor moved to D6)
$ws.Range("D6").Select()
